$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. dataset_type sheet: remove "nanoPOTS" (row 3) and "NanoDESI" (row 21),
#    then add the new "2D Imaging Mass Cytometry" dataset type after "MALDI".
# ---------------------------------------------------------------------------
$dsType = $wb.Sheets.Item("dataset_type")

# Remove row 3 ("nanoPOTS")
$dsType.Rows.Item(3).Delete()

# "NanoDESI" was row 21; after removing row 3 above it is now row 20.
$dsType.Rows.Item(20).Delete()

# "MALDI" is now at row 22 (24 - 2 removed rows above it); insert the new
# dataset type right after it, before "RNAseq (GeoMx)".
$dsType.Rows.Item(23).Insert()
$dsType.Range("A23").Value = "2D Imaging Mass Cytometry"
$dsType.Range("B23").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000296"

# ---------------------------------------------------------------------------
# 2. acquisition_instrument_model sheet: add "STELLARIS 5" after "SCN400"
#    and "Unknown" after "Resolve Biosciences Molecular Cartography".
# ---------------------------------------------------------------------------
$acqModel = $wb.Sheets.Item("acquisition_instrument_model")

# Insert "STELLARIS 5" as the new row 2 (after SCN400, before BZ-X710).
$acqModel.Rows.Item(2).Insert()
$acqModel.Range("A2").Value = "STELLARIS 5"
$acqModel.Range("B2").Value = "https://identifiers.org/RRID:SCR_024663"

# "Resolve Biosciences Molecular Cartography" was row 10; after the insert
# above it is now row 11. Insert "Unknown" right after it, as the new row 12.
$acqModel.Rows.Item(12).Insert()
$acqModel.Range("A12").Value = "Unknown"
$acqModel.Range("B12").Value = "http://ncicb.nci.nih.gov/xml/owl/EVS/Thesaurus.owl#C17998"

# ---------------------------------------------------------------------------
# 3. Update the data-validation list ranges on the SIMS sheet that point at
#    the two lookup sheets above, to reflect their new row counts.
# ---------------------------------------------------------------------------
$sims = $wb.Sheets.Item("SIMS")

$dv = $sims.Range("D2:D1001").Validation
$dv.Modify(3, 1, 1, "'dataset_type'!`$A`$1:`$A`$35")

$dv2 = $sims.Range("H2:H1001").Validation
$dv2.Modify(3, 1, 1, "'acquisition_instrument_model'!`$A`$1:`$A`$38")

# ---------------------------------------------------------------------------
# 4. Bump the metadata "pav:createdOn" timestamp.
# ---------------------------------------------------------------------------
$meta = $wb.Sheets.Item(".metadata")
$meta.Range("C2").Value = "2023-11-02T15:46:47-07:00"
